# Commit: "Wed, Apr 29, 2020  3:06:36 AM"
#
# 1) Slide 16's table (3rd shape) switches its table style from
#    {1526B983-BB0C-4E30-8F27-C087AB898D03} to {D54922CF-8A82-4010-9DF9-6A87EDD5C6EE}.
# 2) The deck's theme colour scheme (used by the single Slide Master / ppt/theme/theme2.xml)
#    changes from the "Integral" palette to the default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{D54922CF-8A82-4010-9DF9-6A87EDD5C6EE}")

# --- 2. Theme colours: Integral -> Office Theme ---------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# index : scheme slot : new RGB (Office Theme palette)
$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
